$d = $word.ActiveDocument

# Replace the city "Market Drayton" with "U" (first run of the location line)
$d.Content.Find.Execute("Market Drayton", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "U", 2)

# Replace the remaining ", UK" with "nited Kingdom" so the line reads "United Kingdom"
$d.Content.Find.Execute(", UK", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "nited Kingdom", 2)
